$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 3 through 13 (old event rows), leaving only header (row1) + 1 data row (row2)
$ws.Rows("3:13").Delete()

# Update the single remaining data row (row 2)
$ws.Range("A2").Value = (Get-Date -Year 2024 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = "Python é melhor que Power Automate e qualquer outra ferraemnte de automação, com o maker!!! 😎"

# Resize column C to fit the new, longer text (matches Excel's bestFit autosize result)
$ws.Columns("C").ColumnWidth = 93.5

# Update the active selection shown in the saved sheet view
$ws.Range("C9").Select()

$wb.Save()
